# renamed exchange_rate to fx and renamed some of the results
$wb = $excel.ActiveWorkbook

# --- workbook / window level changes -------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 21885
$win.Top = 2940
$win.Width = 28800
$win.Height = 15825

# --- add the new sheet at the end ------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Experiments already ran"

# Write cell values in the exact order the strings were first introduced so the
# shared-string table comes out in the same sequence as the target workbook.
$ws3.Range("A1").Value = "Experiments "
$ws3.Range("C3").Value = "TOD"
$ws3.Range("D4").Value = "Index"
$ws3.Range("D21").Value = "Commodity"
$ws3.Range("D26").Value = "FX"
$ws3.Range("D33").Value = "Crypto"
$ws3.Range("E5").Value = "S&P500"
$ws3.Range("E9").Value = "FTSE100"
$ws3.Range("E13").Value = "NASDAQ"
$ws3.Range("E17").Value = "DOWJ"
$ws3.Range("E22").Value = "WTI"
$ws3.Range("E27").Value = "USD/GBP"
$ws3.Range("E34").Value = "BTC"
$ws3.Range("F6").Value = "minutely"
$ws3.Range("F7").Value = "daily"
$ws3.Range("F8").Value = "weekly"

# Remaining cells reuse the strings created above.
$ws3.Range("F10").Value = "minutely"
$ws3.Range("F11").Value = "daily"
$ws3.Range("F12").Value = "weekly"
$ws3.Range("F14").Value = "minutely"
$ws3.Range("F15").Value = "daily"
$ws3.Range("F16").Value = "weekly"
$ws3.Range("F18").Value = "minutely"
$ws3.Range("F19").Value = "daily"
$ws3.Range("F20").Value = "weekly"
$ws3.Range("F23").Value = "daily"
$ws3.Range("F24").Value = "weekly"
$ws3.Range("F28").Value = "minutely"
$ws3.Range("F29").Value = "daily"
$ws3.Range("F30").Value = "weekly"
$ws3.Range("F34").Value = "minutely"
$ws3.Range("F35").Value = "daily"

# --- sheet view / selection state ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("notes and data restrictions")
$ws2.Activate()
$ws2.Range("E11").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

$ws1 = $wb.Worksheets.Item("Experiments")
$ws1.Activate()
$ws1.Range("E25").Select()

$ws3.Activate()
$ws3.Range("J25").Select()
